# Appends new Dice job listing rows (152-158) to the job list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(
        "Golang developer Remote ( But need candidates from Dallas only)",
        "https://www.dice.com/job-detail/a1e3a639-156c-4016-9928-e8c68da6b8c8?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Dallas, Texas",
        "Contract, Third Party",
        "`$60 - `$70",
        "NimbusAITech LLC"
    ),
    @(
        "Go-Lang Developer",
        "https://www.dice.com/job-detail/c017efaa-ec6e-416d-8e2c-e96f495229eb?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Texas",
        "Contract",
        "USD78 - USD77",
        "Pyramid Consulting, Inc."
    ),
    @(
        "Go-Lang Developer",
        "https://www.dice.com/job-detail/f6fb400d-4325-4f52-b9fb-5ebccf178382?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=1&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Hybrid in Westlake, Texas",
        "Contract",
        "74 - 76",
        "SS Software Solutions LLC"
    ),
    @(
        "SAP PP Consultant || Onsite in Normal/Chicago, IL (Relocation will work) || Automotive industry exp.",
        "https://www.dice.com/job-detail/1525e71a-a575-4d66-96a2-d6cc955fcda7?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=3&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Texas",
        "Contract, Third Party",
        "Depends on Experience",
        "TechVirtue LLC"
    ),
    @(
        "Gen AI with Google Cloud Platform",
        "https://www.dice.com/job-detail/443ed1e2-c7f1-42dc-8d40-e82b34cca2e5?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=3&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Remote",
        "Contract",
        "Depends on Experience",
        "Oraapps Inc"
    ),
    @(
        "AWS Full-Stack Backend Developer",
        "https://www.dice.com/job-detail/fefda5a0-978d-41c2-9f3e-1c68614653de?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=4&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Albany, New York",
        "Full-time, Third Party, Contract",
        "Depends on Experience",
        "Calance"
    ),
    @(
        "Backend Java Developer,W2, Need Local Citizens",
        "https://www.dice.com/job-detail/0451a22b-7b3c-4005-b7cb-728fa2eb195d?filters.employmentType=CONTRACTS%7CTHIRD_PARTY&filters.postedDate=ONE&page=4&location=United+States&longitude=-106.5348379&latitude=38.7945952&locationPrecision=Country&countryCode=US&q=Golang",
        "Hybrid in Alpharetta, Georgia",
        "Contract",
        "50 - 60",
        "The Silicon Partners Inc."
    )
)

$startRow = 152
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rowData = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $rowData[0]
    $ws.Cells.Item($r, 2).Value = $rowData[1]
    $ws.Cells.Item($r, 3).Value = $rowData[2]
    $ws.Cells.Item($r, 4).Value = $rowData[3]
    $ws.Cells.Item($r, 5).Value = $rowData[4]
    $ws.Cells.Item($r, 6).Value = $rowData[5]
}
